$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect the merged/changed data
$ws.Range("B2").Value = "CFINHRFLA"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 16960.8644
$ws.Range("F2").Value = "21/03/2025"
$ws.Range("G2").Value = "21/03/2025"
$ws.Range("H2").Value = 33922

# Delete row 3 entirely (shifts cells up, removes the row)
$ws.Range("A3:H3").Delete(-4162) | Out-Null
